$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '27.168.56'
$ws.Range('E2').Value = '  +0.37%  '
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '1.902.27'
$ws.Range('E3').Value = '  +0.83%  '
$ws.Range('E4').Value = '  +0.13%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '306.07'
$ws.Range('E5').Value = '  -0.36%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '0.9999'
$ws.Range('E6').Value = '  +0.11%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.5249'
$ws.Range('E7').Value = '  +1.31%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.3775'
$ws.Range('E8').Value = '  +1.49%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.07244'
$ws.Range('E9').Value = '  +0.57%  '
$ws.Range('E10').Value = '  +0.85%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.8981'
$ws.Range('E11').Value = '  -0.70%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.08367'
$ws.Range('E12').Value = '  +9.62%  '
$ws.Range('E13').Value = '  +1.24%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '94.79'
$ws.Range('E14').Value = '  -0.40%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '5.267'
$ws.Range('E15').Value = '  -0.01%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '1.001'
$ws.Range('E16').Value = '  +0.18%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '0.000008632'
$ws.Range('E17').Value = '  +1.49%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '14.54'
$ws.Range('E18').Value = '  +1.75%  '
$ws.Range('E19').Value = '  +0.08%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '27.204.50'
$ws.Range('E20').Value = '  +0.32%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '5.059'
$ws.Range('E21').Value = '  +0.29%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '2.146.25'
$ws.Range('E22').Value = '  +1.40%  '
$ws.Range('E23').Value = '  +0.79%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '6.425'
$ws.Range('E24').Value = '  -0.57%  '
$ws.Range('B25').Value = 'LidoDAOToken'
$ws.Range('C25').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '2.283'
$ws.Range('E25').Value = '  +7.56%  '
$ws.Range('B26').Value = 'Monero'
$ws.Range('C26').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '146.53'
$ws.Range('E26').Value = '  +0.57%  '
$ws.Range('E27').Value = '  -1.82%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '18.12'
$ws.Range('E28').Value = '  +0.54%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '114.80'
$ws.Range('E29').Value = '  +0.04%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '4.926'
$ws.Range('E30').Value = '  +0.01%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '4.786'
$ws.Range('E31').Value = '  -0.17%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '0.09264'
$ws.Range('E32').Value = '  +0.68%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '0.8206'
$ws.Range('E33').Value = '  +7.58%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '0.05055'
$ws.Range('E34').Value = '  +0.31%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '1.235'
$ws.Range('E35').Value = '  +3.80%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '2.958'
$ws.Range('E36').Value = '  -1.94%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '3.349'
$ws.Range('E37').Value = '  +2.13%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '2.603'
$ws.Range('E38').Value = '  +2.67%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.5709'
$ws.Range('E39').Value = '  +1.73%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.01981'
$ws.Range('E40').Value = '  -0.54%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '1.071'
$ws.Range('E41').Value = '  -0.52%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '6.678'
$ws.Range('E42').Value = '  +1.26%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '8.954'
$ws.Range('E43').Value = '  +1.14%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '118.00'
$ws.Range('E44').Value = '  -0.05%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '0.1512'
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '0.4840'
$ws.Range('E46').Value = '  +1.02%  '
$ws.Range('B47').Value = 'EnergySwap'
$ws.Range('C47').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '10.16'
$ws.Range('E47').Value = '  +0.25%  '
$ws.Range('B48').Value = 'PaxDollar'
$ws.Range('C48').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '0.9997'
$ws.Range('E48').Value = '  +0.11%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '1.615'
$ws.Range('E49').Value = '  +2.67%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '37.43'
$ws.Range('E50').Value = '  +0.81%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '63.65'
$ws.Range('E51').Value = '  +0.19%  '
